$wb = $excel.ActiveWorkbook
$general = $wb.Worksheets.Item("General")
$indexes = $wb.Worksheets.Item("Indexes.xlsx")

# ---- General sheet numeric updates ----
$general.Range("H2").Value = 120.0
$general.Range("I2").Value = 110.0
$general.Range("J2").Value = 100.0
$general.Range("M2").Value = 65.0
$general.Range("N2").Value = 65.0
$general.Range("O2").Value = 60.0
$general.Range("R2").Value = "-"
$general.Range("T2").Value = "-"

$general.Range("H3").Value = 6.666666666666666
$general.Range("I3").Value = 6.11
$general.Range("J3").Value = 5.555555555555555
$general.Range("M3").Value = 390.0
$general.Range("N3").Value = 390.0
$general.Range("O3").Value = 360.0
$general.Range("R3").Value = "-"
$general.Range("T3").Value = "-"

# ---- Indexes.xlsx sheet updates ----
$indexes.Range("C2").Value = "belfiore"
$indexes.Range("D2").Formula = "=2 / (((0.5 * General!F3 + General!H3 + General!J3) / 19.08) * ((0.5 * General!L2 + General!N2 + General!P2) / 104.0) + 1 )"
$indexes.Range("E2").Value = "Healthy"
$indexes.Range("F2").Value = "≅1"

$indexes.Range("C4").Value = "avingon"
$indexes.Range("D4").Formula = "=((0.137 * 100000000 /(General!F3 * General!L2 * 150/General!Q3)) + 100000000 /(General!J3 * General!P2 * 150/General!Q3)) / 2"
$indexes.Range("F4").Value = "-"
